$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update IPC PO (C), DELTA (D), DELTA^2 (E) for rows 2-51
$cData = @(
  29.72175776911019,
  30.02891384579398,
  29.87532757460375,
  29.64943781981881,
  30.2197193266033,
  30.16524852836157,
  30.95421275670812,
  30.3238393500008,
  29.46879379716013,
  30.68162427992516,
  30.02526259058906,
  29.65885848788355,
  30.0738608790388,
  29.9260358867717,
  29.94177716779699,
  29.90366860471718,
  29.86792644700605,
  30.78666544642995,
  30.99431524769091,
  31.60933116519394,
  31.35647702173609,
  31.14129646547685,
  31.11198262567547,
  31.61760285234385,
  31.5746255466973,
  31.46276909702048,
  31.6124046137601,
  32.07227687734462,
  32.43597484441806,
  32.75777595936748,
  34.12422113804068,
  33.96927027946951,
  33.53024531804003,
  32.84650728479531,
  34.00686701367583,
  34.89016021427642,
  35.13709982289468,
  35.28494508555107,
  34.84355280154394,
  35.66862906993612,
  35.85962750787648,
  36.58988521355371,
  37.83083254539498,
  39.01314395352769,
  39.25621049148977,
  39.2911214137161,
  40.10510052842852,
  39.87448722716134,
  36.32202583591969,
  39.63872335809702
)
$dData = @(
  0.7417577691101904,
  0.8789138457939778,
  0.5253275746037467,
  0.2794378198188099,
  0.6797193266032977,
  0.6152485283615654,
  1.20421275670812,
  0.483839350000796,
  -0.3412062028398708,
  0.761624279925158,
  0.04526259058906135,
  -0.3811415121164501,
  -0.1361391209612037,
  -0.2939641132282951,
  -0.4382228322030066,
  -0.5363313952828221,
  -0.6120735529939552,
  0.09666544642994523,
  0.2443152476909134,
  0.6693311651939347,
  0.4064770217360909,
  0.1212964654768527,
  -0.008017374324534643,
  0.3376028523438528,
  0.1946255466973028,
  -0.117230902979518,
  -0.03759538623989656,
  0.1922768773446215,
  0.1559748444180542,
  0.3077759593674756,
  1.27422113804068,
  1.069270279469514,
  0.4302453180400292,
  -0.5534927152046905,
  0.3068670136758271,
  0.790160214276419,
  0.7370998228946846,
  0.3849450855510739,
  -0.4564471984560612,
  -0.03137093006388625,
  -0.4403724921235153,
  -0.2101147864462831,
  0.5308325453949863,
  1.11314395352769,
  0.7562104914897674,
  0.3911214137160997,
  0.7051005284285239,
  -0.02551277283865971,
  -3.777974164080312,
  -0.9612766419029768
)
$eData = @(
  0.5502045880353266,
  0.7724895483283601,
  0.275969060639055,
  0.07808549514508965,
  0.4620183629580405,
  0.378530751651072,
  1.45012836341857,
  0.2341005166091928,
  0.116421672856403,
  0.5800715437715155,
  0.002048702106832985,
  0.1452688522584141,
  0.01853386025608926,
  0.08641489986609789,
  0.1920392506640245,
  0.2876513655660188,
  0.374634034274644,
  0.009344208533500611,
  0.05968994025427236,
  0.4480042086998702,
  0.1652235691994425,
  0.01471283253717732,
  0.00006427829105970732,
  0.1139756859107053,
  0.037879103427224,
  0.01374308461339316,
  0.001413413066527003,
  0.03697039756139862,
  0.02432815209123622,
  0.09472604116457,
  1.623639508629686,
  1.143338930556812,
  0.1851110336953659,
  0.3063541857846606,
  0.09416736408232028,
  0.6243531642253565,
  0.5433161489113754,
  0.1481827188899236,
  0.2083440449783869,
  0.0009841352530732422,
  0.1939279318190755,
  0.04414822348336716,
  0.2817831912505201,
  1.239089461275255,
  0.5718543074391955,
  0.1529759602672804,
  0.4971667551901836,
  0.0006509015779170528,
  14.27308878445833,
  0.9240527822682639
)

for ($i = 0; $i -lt $cData.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 3).Value = $cData[$i]
  $ws.Cells.Item($row, 4).Value = $dData[$i]
  $ws.Cells.Item($row, 5).Value = $eData[$i]
}

$ws.Cells.Item(52, 3).Value = 8.072418978433124
$ws.Cells.Item(52, 5).Value = 30.08121531779148
$ws.Cells.Item(53, 5).Value = 0.6016243063558295

Write-Host "Done updating sliding window results"